$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing C-column values (rows 2-64) ---
$ws.Range("C2").Value = 0.015622
$ws.Range("C3").Value = 0.031255
$ws.Range("C4").Value = 0.004948
$ws.Range("C5").Value = 0.038367
$ws.Range("C6").Value = 0
$ws.Range("C19").Value = 0.013658
$ws.Range("C20").Value = 0.000509
$ws.Range("C24").Value = 0
$ws.Range("C32").Value = 0.013834
$ws.Range("C37").Value = 0.006093
$ws.Range("C38").Value = 0.000607
$ws.Range("C48").Value = 0.015644
$ws.Range("C49").Value = 0.01568
$ws.Range("C50").Value = 0.015622
$ws.Range("C51").Value = 0.015633
$ws.Range("C52").Value = 0.015632
$ws.Range("C53").Value = 0.01561
$ws.Range("C54").Value = 0
$ws.Range("C55").Value = 0.017542
$ws.Range("C56").Value = 0.025542
$ws.Range("C57").Value = 0.015014
$ws.Range("C58").Value = 0
$ws.Range("C59").Value = 0.015707
$ws.Range("C60").Value = 0.015629
$ws.Range("C61").Value = 0.015622
$ws.Range("C62").Value = 0.015624
$ws.Range("C63").Value = 0.015623
$ws.Range("C64").Value = 0.01563

# --- Append new rows 65-84 (white_fred_1.jpg .. white_fred_20.jpg) ---
$ws.Range("B65:B84").NumberFormat = "@"
$ws.Range("A65").Value = "white_fred_1.jpg"
$ws.Range("B65").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C65").Value = 0.015573
$ws.Range("D65").Value = 0
$ws.Range("A66").Value = "white_fred_2.jpg"
$ws.Range("B66").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C66").Value = 0.022932
$ws.Range("D66").Value = 0
$ws.Range("A67").Value = "white_fred_3.jpg"
$ws.Range("B67").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C67").Value = 0.005103
$ws.Range("D67").Value = 0
$ws.Range("A68").Value = "white_fred_4.jpg"
$ws.Range("B68").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C68").Value = 0.015696
$ws.Range("D68").Value = 0
$ws.Range("A69").Value = "white_fred_5.jpg"
$ws.Range("B69").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C69").Value = 0.015637
$ws.Range("D69").Value = 0
$ws.Range("A70").Value = "white_fred_6.jpg"
$ws.Range("B70").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C70").Value = 0.015613
$ws.Range("D70").Value = 0
$ws.Range("A71").Value = "white_fred_7.jpg"
$ws.Range("B71").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C71").Value = 0.015624
$ws.Range("D71").Value = 0
$ws.Range("A72").Value = "white_fred_8.jpg"
$ws.Range("B72").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C72").Value = 0.01562
$ws.Range("D72").Value = 0
$ws.Range("A73").Value = "white_fred_9.jpg"
$ws.Range("B73").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C73").Value = 0.015625
$ws.Range("D73").Value = 0
$ws.Range("A74").Value = "white_fred_10.jpg"
$ws.Range("B74").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C74").Value = 0.015632
$ws.Range("D74").Value = 0
$ws.Range("A75").Value = "white_fred_11.jpg"
$ws.Range("B75").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("A76").Value = "white_fred_12.jpg"
$ws.Range("B76").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C76").Value = 0.015629
$ws.Range("D76").Value = 0
$ws.Range("A77").Value = "white_fred_13.jpg"
$ws.Range("B77").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C77").Value = 0.015622
$ws.Range("D77").Value = 0
$ws.Range("A78").Value = "white_fred_14.jpg"
$ws.Range("B78").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C78").Value = 0.015629
$ws.Range("D78").Value = 0
$ws.Range("A79").Value = "white_fred_15.jpg"
$ws.Range("B79").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C79").Value = 0.015618
$ws.Range("D79").Value = 0
$ws.Range("A80").Value = "white_fred_16.jpg"
$ws.Range("B80").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C80").Value = 0.015633
$ws.Range("D80").Value = 0
$ws.Range("A81").Value = "white_fred_17.jpg"
$ws.Range("B81").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C81").Value = 0.015618
$ws.Range("D81").Value = 0
$ws.Range("A82").Value = "white_fred_18.jpg"
$ws.Range("B82").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C82").Value = 0
$ws.Range("D82").Value = 0
$ws.Range("A83").Value = "white_fred_19.jpg"
$ws.Range("B83").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C83").Value = 0.01563
$ws.Range("D83").Value = 0
$ws.Range("A84").Value = "white_fred_20.jpg"
$ws.Range("B84").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C84").Value = 0.015625
$ws.Range("D84").Value = 0
